# Finished new absolute reference tutorial
#
# This script fills in the "Times Table" worksheet (sheet2) with an
# absolute/relative reference multiplication-table tutorial, applies the
# "Accent2"/"Accent5" cell styles to its header cells, and makes the
# "Times Table" sheet the active sheet/tab (it was previously "Payroll").

$wb = $excel.ActiveWorkbook

$wsPayroll = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# Header row: "X" label plus the multipliers 1..10 across B1:K1
# ------------------------------------------------------------------
$ws.Range("A1").Value = "X"

for ($col = 2; $col -le 11; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# ------------------------------------------------------------------
# Column A: the multiplicand 1..10 down rows 2..11
# ------------------------------------------------------------------
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1
}

# ------------------------------------------------------------------
# Formulas:
#   B2 is a standalone (non-shared) formula: =B$1*$A2
#   B3:B11 is one shared-formula group:       =B$1*$A<row>
#   C2:K11 is another shared-formula group:    =<col>$1*$A<row>
# ------------------------------------------------------------------
$ws.Range("B2").Formula = '=B$1*$A2'
$ws.Range("B3:B11").Formula = '=B$1*$A3'
$ws.Range("C2:K11").Formula = '=C$1*$A2'

# ------------------------------------------------------------------
# Styling: label cell gets "Accent5", the header numbers (row 1 and
# column A) get "Accent2" (applied in this order so the style/fill/
# font bookkeeping matches the order they were introduced).
# ------------------------------------------------------------------
$ws.Range("B1:K1").Style = "Accent2"
$ws.Range("A2:A11").Style = "Accent2"
$ws.Range("A1").Style = "Accent5"

# ------------------------------------------------------------------
# Selection / activation: "Times Table" becomes the active sheet with
# K11 selected; "Payroll" loses its tabSelected flag.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("K11").Select()

Write-Output "Times Table tutorial filled in"
